# Updates the cryptocurrency price/volume table to match the latest scrape.
# D-column "Price" values and E-column "Volume(1h)" percentages are stored
# as plain text in this workbook (not numbers), so values that otherwise
# look numeric (e.g. "1.010", "347.14") must be forced to remain text —
# Excel would otherwise silently reinterpret them as numbers and drop
# trailing zeros / collapse formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is unambiguously non-numeric (contains letters,
# a percent sign, multiple dots, surrounding spaces, etc.) - safe to set
# directly without touching number format.
$textUpdates = @{
    "D2" = "29.965.42"
    "E2" = "  -0.03%  "
    "D3" = "2.116.66"
    "E3" = "  +0.79%  "
    "E4" = "  +0.57%  "
    "E5" = "  +0.21%  "
    "E6" = "  +0.59%  "
    "E7" = "  +0.78%  "
    "E8" = "  +0.47%  "
    "E9" = "  +2.19%  "
    "E10" = "  -0.07%  "
    "E11" = "  +0.82%  "
    "E12" = "  -0.43%  "
    "E13" = "  +4.68%  "
    "D14" = "2.132.58"
    "E14" = "  +1.73%  "
    "E15" = "  +2.29%  "
    "E16" = "  +3.34%  "
    "E17" = "  +0.19%  "
    "E18" = "  +0.56%  "
    "E19" = "  +4.54%  "
    "E20" = "  +0.38%  "
    "E21" = "  +1.41%  "
    "E22" = "  +0.57%  "
    "D23" = "30.018.11"
    "E23" = "  -0.16%  "
    "E24" = "  +0.43%  "
    "E25" = "  -0.30%  "
    "D26" = "2.380.83"
    "E26" = "  +1.60%  "
    "E27" = "  +0.24%  "
    "E28" = "  -0.59%  "
    "E29" = "  -0.10%  "
    "E30" = "  +0.54%  "
    "E31" = "  -1.62%  "
    "E32" = "  +7.81%  "
    "E33" = "  -0.20%  "
    "E34" = "  +0.13%  "
    "B35" = "InternetComputer(DFINITY)"
    "C35" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "E35" = "  +5.93%  "
    "B36" = "HuobiToken"
    "C36" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "E36" = "  +0.79%  "
    "E37" = "  +5.48%  "
    "E38" = "  +2.38%  "
    "E39" = "  +1.49%  "
    "E40" = "  +1.90%  "
    "E41" = "  +1.15%  "
    "E42" = "  -1.65%  "
    "E43" = "  +1.80%  "
    "E44" = "  +2.98%  "
    "E45" = "  +2.52%  "
    "E46" = "  +3.36%  "
    "E47" = "  +0.57%  "
    "E48" = "  +3.35%  "
    "E49" = "  +10.01%  "
    "E50" = "  +0.05%  "
    "E51" = "  +0.26%  "
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Cells whose new text would otherwise be auto-coerced into a numeric
# value by Excel (single decimal point, no other characters) - force the
# cell to text format first so the literal string is preserved exactly.
$numericLookingUpdates = @{
    "D4" = "1.010"
    "D5" = "347.14"
    "D7" = "0.5192"
    "D8" = "0.4446"
    "D9" = "53.65"
    "D10" = "0.09331"
    "D12" = "25.14"
    "D13" = "8.540"
    "D15" = "6.898"
    "D16" = "102.76"
    "D19" = "21.51"
    "D20" = "0.06702"
    "D21" = "6.306"
    "D22" = "1.008"
    "D24" = "12.69"
    "D25" = "2.324"
    "D27" = "22.07"
    "D28" = "2.537"
    "D29" = "162.53"
    "D30" = "134.14"
    "D31" = "1.147"
    "D32" = "1.768"
    "D34" = "6.242"
    "D35" = "6.578"
    "D36" = "3.977"
    "D37" = "10.72"
    "D38" = "0.02623"
    "D39" = "0.06879"
    "D40" = "0.7054"
    "D41" = "12.67"
    "D42" = "0.2242"
    "D43" = "1.330"
    "D44" = "0.6842"
    "D45" = "14.52"
    "D46" = "2.356"
    "D47" = "1.008"
    "D48" = "0.00000000363"
    "D49" = "1.276"
}

foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
}
